$wb = $excel.ActiveWorkbook

# --- Update the "Last update" timestamp on the "info" sheet ---
$infoSheet = $wb.Worksheets.Item("info")
$infoSheet.Range("B2").Value = "2021-10-15 15:19:26"

# --- Append United Kingdom and France rows to the "panel b-e" sheet ---
$sheet = $wb.Worksheets.Item("panel b-e")

# Row 20: United Kingdom
$sheet.Cells.Item(20, 1).Value = "United Kingdom"
$sheet.Cells.Item(20, 2).Value = "GBR"
$sheet.Cells.Item(20, 3).Value = "DEV"
$sheet.Cells.Item(20, 4).Value = -3.13681919422224
$sheet.Cells.Item(20, 5).Value = -0.149748265115072
$sheet.Cells.Item(20, 6).Value = 6.74437083964458
$sheet.Cells.Item(20, 7).Value = 0.145332580540567

# Row 21: France
$sheet.Cells.Item(21, 1).Value = "France"
$sheet.Cells.Item(21, 2).Value = "FRA"
$sheet.Cells.Item(21, 3).Value = "DEV"
$sheet.Cells.Item(21, 4).Value = -1.62051639729969
$sheet.Cells.Item(21, 5).Value = -0.0707187356488979
$sheet.Cells.Item(21, 6).Value = 6.63879934291407
$sheet.Cells.Item(21, 7).Value = 0.144843896816565
